$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.411.18'
$ws.Range('E2').Value = '  -0.05%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.573.47'
$ws.Range('E3').Value = '  +0.16%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.16%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.41'
$ws.Range('E6').Value = '  -0.50%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3769'
$ws.Range('E7').Value = '  +3.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.90'
$ws.Range('E8').Value = '  +1.16%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3425'
$ws.Range('E9').Value = '  +1.16%  '

$ws.Range('E10').Value = '  -1.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07653'
$ws.Range('E11').Value = '  +0.89%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.15%  '

$ws.Range('E13').Value = '  +0.50%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.026'
$ws.Range('E14').Value = '  -0.70%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.950'
$ws.Range('E15').Value = '  +0.86%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.573.09'
$ws.Range('E16').Value = '  +0.04%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001135'
$ws.Range('E17').Value = '  -0.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.25'
$ws.Range('E18').Value = '  +1.24%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06761'
$ws.Range('E19').Value = '  +0.45%  '

$ws.Range('E20').Value = '  +0.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.88'
$ws.Range('E21').Value = '  +2.59%  '

$ws.Range('E22').Value = '  -0.50%  '

$ws.Range('E23').Value = '  -0.20%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.429'
$ws.Range('E24').Value = '  +1.56%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '22.403.30'
$ws.Range('E25').Value = '  -0.10%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.719'
$ws.Range('E26').Value = '  -10.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.32'
$ws.Range('E27').Value = '  +2.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '147.28'
$ws.Range('E28').Value = '  +2.05%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.033'
$ws.Range('E29').Value = '  +1.28%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.38'
$ws.Range('E30').Value = '  +0.84%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.750.19'
$ws.Range('E31').Value = '  +0.02%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.192'
$ws.Range('E32').Value = '  -1.37%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.011'
$ws.Range('E33').Value = '  +1.39%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9938'
$ws.Range('E34').Value = '  -6.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.04'
$ws.Range('E35').Value = '  -3.39%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08599'
$ws.Range('E36').Value = '  +1.74%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02552'
$ws.Range('E37').Value = '  -0.51%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2318'
$ws.Range('E38').Value = '  +0.59%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06581'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.332'
$ws.Range('E40').Value = '  +6.69%  '

$ws.Range('E41').Value = '  -1.17%  '

$ws.Range('E42').Value = '  +0.94%  '

$ws.Range('E43').Value = '  -2.66%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.14'
$ws.Range('E44').Value = '  -3.19%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.002'
$ws.Range('E45').Value = '  +0.19%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.795'
$ws.Range('E46').Value = '  +0.36%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6016'
$ws.Range('E47').Value = '  -0.24%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.306'
$ws.Range('E48').Value = '  +7.51%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.087'
$ws.Range('E49').Value = '  -2.17%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.70'
$ws.Range('E50').Value = '  +1.80%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07330'
$ws.Range('E51').Value = '  +0.44%  '
